$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches existing text-cell semantics)
$textGuardCells = @("D5", "D6", "D11", "D12", "D15", "D16", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D37", "D38", "D39", "D40", "D43", "D47", "D50", "D51")
foreach ($addr in $textGuardCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.167.76'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '3.317.61'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '564.70'
$ws.Range("E5").Value = '  +1.55%  '
$ws.Range("D6").Value = '186.09'
$ws.Range("E6").Value = '  +2.22%  '
$ws.Range("D8").Value = '3.311.36'
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("E10").Value = '  -4.21%  '
$ws.Range("D11").Value = '0.573'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").Value = '46.16'
$ws.Range("E12").Value = '  -2.12%  '
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '3.847.77'
$ws.Range("E14").Value = '  +0.79%  '
$ws.Range("D15").Value = '8.45'
$ws.Range("E15").Value = '  -1.98%  '
$ws.Range("D16").Value = '595.78'
$ws.Range("E16").Value = '  -7.62%  '
$ws.Range("D17").Value = '66.122.13'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").Value = '3.316.51'
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("D20").Value = '17.71'
$ws.Range("E20").Value = '  -1.95%  '
$ws.Range("D21").Value = '10.89'
$ws.Range("E21").Value = '  -4.06%  '
$ws.Range("D22").Value = '0.896'
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").Value = '17.90'
$ws.Range("E23").Value = '  -0.94%  '
$ws.Range("D24").Value = '5.00'
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("D25").Value = '98.88'
$ws.Range("E25").Value = '  -7.90%  '
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = '2.72'
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("D28").Value = '9.42'
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '30.77'
$ws.Range("E29").Value = '  +1.98%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '8.45'
$ws.Range("E30").Value = '  -2.26%  '
$ws.Range("D31").Value = '6.64'
$ws.Range("E31").Value = '  +6.18%  '
$ws.Range("E32").Value = '  -5.83%  '
$ws.Range("D33").Value = '559.46'
$ws.Range("E33").Value = '  +7.28%  '
$ws.Range("D34").Value = '10.86'
$ws.Range("E34").Value = '  -1.32%  '
$ws.Range("D35").Value = '3.804.19'
$ws.Range("E35").Value = '  +0.67%  '
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").Value = '55.93'
$ws.Range("E38").Value = '  -2.40%  '
$ws.Range("D39").Value = '33.24'
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").Value = '0.127'
$ws.Range("E40").Value = '  -1.04%  '
$ws.Range("D41").Value = '0.0₃0686'
$ws.Range("E41").Value = '  -5.86%  '
$ws.Range("E42").Value = '  -5.67%  '
$ws.Range("D43").Value = '3.39'
$ws.Range("E43").Value = '  +3.97%  '
$ws.Range("E44").Value = '  -3.67%  '
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("D47").Value = '3.05'
$ws.Range("E47").Value = '  -8.19%  '
$ws.Range("E48").Value = '  -2.08%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").Value = '2.52'
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("D51").Value = '129.46'
$ws.Range("E51").Value = '  +6.09%  '
